# Update the "Quy trinh thuc hien" (process of execution) slide body text.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

# --- Buoc 1 : just reword the single run in place ---
$para1 = $tr.Paragraphs(1)
$para1.Runs(1).Text = "Bước 1: thu thập dữ liệu ảnh."

# --- Buoc 2 : collapse the three runs into a single run with new wording ---
$para2 = $tr.Paragraphs(2)
# remove the extra runs first (highest index first so indices stay valid)
$para2.Runs(3).Text = ""
$para2.Runs(2).Text = ""
$para2.Runs(1).Text = "Bước 2: gán nhãn thủ công cho dữ liệu đã thu thập bằng roboflow."
# the explicit 115% line spacing override on this paragraph was dropped
$para2.ParagraphFormat.SpaceWithin = 1.0

# --- Buoc 3 : reword the single run in place ---
$para3 = $tr.Paragraphs(3)
$para3.Runs(1).Text = "Bước 3: sử dụng trang roboflow tăng cường thêm dữ liệu trên và trích xuất file TF-Record."
# the explicit 115% line spacing override on this paragraph was dropped too
$para3.ParagraphFormat.SpaceWithin = 1.0
